# Auto-generated cell updates applying the Cactuar_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 880.03125
$ws.Range("I28").Value = 363.86667
$ws.Range("K28").Value = 363.86667
$ws.Range("M28").Value = 121.13333
$ws.Range("H51").Value = 5406
$ws.Range("I51").Value = 5408.3335
$ws.Range("J51").Value = 5403.6665
$ws.Range("K51").Value = 5408.3335
$ws.Range("L51").Value = 5403.6665
$ws.Range("M51").Value = -4924.3335
$ws.Range("N51").Value = -6371.6665
$ws.Range("H62").Value = 3014.1875
$ws.Range("I62").Value = 2483
$ws.Range("K62").Value = 2483
$ws.Range("M62").Value = -1859
$ws.Range("H65").Value = 3014.1875
$ws.Range("I65").Value = 2483
$ws.Range("K65").Value = 12415
$ws.Range("M65").Value = -9295
$ws.Range("H70").Value = 6500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 19500
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -20040
$ws.Range("H73").Value = 6500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 19500
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -21372
$ws.Range("H98").Value = 906.2308
$ws.Range("I98").Value = 815.0833
$ws.Range("K98").Value = 815.0833
$ws.Range("M98").Value = 682.9167
$ws.Range("H111").Value = 2147.75
$ws.Range("I111").Value = 2164
$ws.Range("K111").Value = 6492
$ws.Range("M111").Value = -3425
$ws.Range("H113").Value = 44380.41
$ws.Range("I113").Value = 2948.5
$ws.Range("J113").Value = 49904.668
$ws.Range("K113").Value = 2948.5
$ws.Range("L113").Value = 49904.668
$ws.Range("M113").Value = 305.5
$ws.Range("N113").Value = -56412.668
$ws.Range("H121").Value = 3913.1667
$ws.Range("J121").Value = 3913.1667
$ws.Range("L121").Value = 11739.5001
$ws.Range("N121").Value = -15233.5001
$ws.Range("H122").Value = 906.2308
$ws.Range("I122").Value = 815.0833
$ws.Range("K122").Value = 2445.2499
$ws.Range("M122").Value = 4.750100000000202
$ws.Range("H131").Value = 3790.5715
$ws.Range("I131").Value = 3272.25
$ws.Range("J131").Value = 4481.6665
$ws.Range("K131").Value = 9816.75
$ws.Range("L131").Value = 13444.9995
$ws.Range("M131").Value = -4776.75
$ws.Range("N131").Value = -23524.9995
$ws.Range("H132").Value = 134802.34
$ws.Range("I132").Value = 257024.28
$ws.Range("K132").Value = 771072.84
$ws.Range("M132").Value = -768542.84
$ws.Range("H138").Value = 3221.5
$ws.Range("I138").Value = 1540.75
$ws.Range("K138").Value = 4622.25
$ws.Range("M138").Value = 517.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4319.6665
$ws.Range("I32").Value = 2153.0417
$ws.Range("K32").Value = 2153.0417
$ws.Range("M32").Value = -1866.0417
$ws.Range("H45").Value = 4749
$ws.Range("I45").Value = 3999
$ws.Range("K45").Value = 3999
$ws.Range("M45").Value = -3622
$ws.Range("H61").Value = 18724
$ws.Range("I61").Value = 18724
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 18724
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -18512
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 25002350
$ws.Range("I74").Value = 41668170
$ws.Range("J74").Value = 3624.5
$ws.Range("K74").Value = 41668170
$ws.Range("L74").Value = 3624.5
$ws.Range("M74").Value = -41667296
$ws.Range("N74").Value = -5372.5
$ws.Range("H77").Value = 25002350
$ws.Range("I77").Value = 41668170
$ws.Range("J77").Value = 3624.5
$ws.Range("K77").Value = 208340850
$ws.Range("L77").Value = 18122.5
$ws.Range("M77").Value = -208336482
$ws.Range("N77").Value = -26858.5
$ws.Range("H82").Value = 119999
$ws.Range("J82").Value = 119999
$ws.Range("L82").Value = 119999
$ws.Range("N82").Value = -120721
$ws.Range("H85").Value = 119999
$ws.Range("J85").Value = 119999
$ws.Range("L85").Value = 119999
$ws.Range("N85").Value = -122495
$ws.Range("H110").Value = 6129.875
$ws.Range("I110").Value = 4008.4
$ws.Range("K110").Value = 4008.4
$ws.Range("M110").Value = -1963.4
$ws.Range("H132").Value = 6144.433
$ws.Range("I132").Value = 2473.9285
$ws.Range("J132").Value = 9356.125
$ws.Range("K132").Value = 7421.7855
$ws.Range("L132").Value = 28068.375
$ws.Range("M132").Value = -4891.7855
$ws.Range("N132").Value = -33128.375
$ws.Range("H136").Value = 18724
$ws.Range("I136").Value = 18724
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 56172
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -53622
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449.75
$ws.Range("I22").Value = 449.75
$ws.Range("K22").Value = 449.75
$ws.Range("M22").Value = -276.75
$ws.Range("H134").Value = 6181.636
$ws.Range("I134").Value = 2255.75
$ws.Range("K134").Value = 6767.25
$ws.Range("M134").Value = -4232.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224
$ws.Range("H22").Value = 792.25
$ws.Range("J22").Value = 771.7143
$ws.Range("L22").Value = 771.7143
$ws.Range("N22").Value = -1471.7143
$ws.Range("H31").Value = 20839384
$ws.Range("I31").Value = 111112920
$ws.Range("J31").Value = 7028.6665
$ws.Range("K31").Value = 111112920
$ws.Range("L31").Value = 7028.6665
$ws.Range("M31").Value = -111112625
$ws.Range("N31").Value = -7618.6665
$ws.Range("H34").Value = 20839384
$ws.Range("I34").Value = 111112920
$ws.Range("J34").Value = 7028.6665
$ws.Range("K34").Value = 111112920
$ws.Range("L34").Value = 7028.6665
$ws.Range("M34").Value = -111112718
$ws.Range("N34").Value = -7432.6665
$ws.Range("H58").Value = 273046.66
$ws.Range("I58").Value = 500683.4
$ws.Range("J58").Value = 5238.706
$ws.Range("K58").Value = 500683.4
$ws.Range("L58").Value = 5238.706
$ws.Range("M58").Value = -500480.4
$ws.Range("N58").Value = -5644.706
$ws.Range("H94").Value = 4202.0625
$ws.Range("I94").Value = 3400.8333
$ws.Range("K94").Value = 3400.8333
$ws.Range("M94").Value = -2949.8333
$ws.Range("H110").Value = 93000
$ws.Range("J110").Value = 93000
$ws.Range("L110").Value = 93000
$ws.Range("N110").Value = -101180
$ws.Range("H136").Value = 273046.66
$ws.Range("I136").Value = 500683.4
$ws.Range("J136").Value = 5238.706
$ws.Range("K136").Value = 1502050.2
$ws.Range("L136").Value = 15716.118
$ws.Range("M136").Value = -1499500.2
$ws.Range("N136").Value = -20816.118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 906.4545000000001
$ws.Range("J7").Value = 1715.6
$ws.Range("L7").Value = 5146.799999999999
$ws.Range("N7").Value = -5370.799999999999
$ws.Range("H99").Value = 4314.8
$ws.Range("J99").Value = 5833
$ws.Range("L99").Value = 17499
$ws.Range("N99").Value = -21991
$ws.Range("H131").Value = 8477172
$ws.Range("I131").Value = 15153673
$ws.Range("J131").Value = 6947140.5
$ws.Range("K131").Value = 45461019
$ws.Range("L131").Value = 20841421.5
$ws.Range("M131").Value = -45455979
$ws.Range("N131").Value = -20851501.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3969041.8
$ws.Range("I107").Value = 7937200
$ws.Range("J107").Value = 883.1667
$ws.Range("K107").Value = 7937200
$ws.Range("L107").Value = 883.1667
$ws.Range("M107").Value = -7935280
$ws.Range("N107").Value = -4723.1667
$ws.Range("H122").Value = 300378
$ws.Range("I122").Value = 369099.84
$ws.Range("J122").Value = 5855.857
$ws.Range("K122").Value = 1107299.52
$ws.Range("L122").Value = 17567.571
$ws.Range("M122").Value = -1104849.52
$ws.Range("N122").Value = -22467.571
$ws.Range("H132").Value = 80912.734
$ws.Range("I132").Value = 136741.33
$ws.Range("J132").Value = 4782.8184
$ws.Range("K132").Value = 410223.99
$ws.Range("L132").Value = 14348.4552
$ws.Range("M132").Value = -407693.99
$ws.Range("N132").Value = -19408.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5490
$ws.Range("I7").Value = 2400
$ws.Range("K7").Value = 2400
$ws.Range("M7").Value = -2288
$ws.Range("H36").Value = 100000
$ws.Range("J36").Value = 100000
$ws.Range("L36").Value = 100000
$ws.Range("N36").Value = -101124
$ws.Range("H126").Value = 5490
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730
$ws.Range("H136").Value = 4106.148
$ws.Range("I136").Value = 2284.1177
$ws.Range("J136").Value = 7203.6
$ws.Range("K136").Value = 6852.353099999999
$ws.Range("L136").Value = 21610.8
$ws.Range("M136").Value = -4302.353099999999
$ws.Range("N136").Value = -26710.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 40000
$ws.Range("I16").Value = 40000
$ws.Range("K16").Value = 40000
$ws.Range("M16").Value = -39708
$ws.Range("H50").Value = 3500
$ws.Range("J50").Value = 3500
$ws.Range("L50").Value = 3500
$ws.Range("N50").Value = -4762
$ws.Range("H55").Value = 26026.5
$ws.Range("J55").Value = 26026.5
$ws.Range("L55").Value = 26026.5
$ws.Range("N55").Value = -26580.5
$ws.Range("H133").Value = 60801.668
$ws.Range("J133").Value = 60801.668
$ws.Range("L133").Value = 60801.668
$ws.Range("N133").Value = -70921.66800000001
